# InformacjeOPrzeniesieniach.xlsx — add a new transfer entry
#
# The "Oddziały" sheet (a table of lesson transfers) gets one new row
# inserted at row 10 (pushing the existing rows 10-15 down to 11-16),
# plus a small text correction in what becomes row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 10.
$ws.Rows.Item(10).Insert()

# Fill in the new row 10 with the new transfer entry.
$ws.Range("A10").Value = "17.12.2025, 7, 13:15-14:00, sala: 41"
$ws.Range("B10").Value = "17.12.2025, 7, 13:15-14:00, sala: 22"
$ws.Range("C10").Value = "Wójcik Kamil"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = "1WA"
$ws.Range("F10").Value = "Zajęcia z wychowawcą"
$ws.Range("G10").Value = ""

# Room-number correction on the row that is now row 16 (was row 15):
# "sala: 40" -> "sala: 4" in the "Przeniesiono na" column.
$ws.Range("B16").Value = "19.12.2025, 9, 14:55-15:40, sala: 4"
